$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.118.82'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '1.835.06'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6835'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3014'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07465'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07659'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.88%  '

$ws.Range("D12").Value = '1.833.03'
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.058'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6806'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.232'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.85%  '

$ws.Range("D17").Value = '29.121.07'
$ws.Range("E17").Value = '  +0.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008192'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").Value = '2.082.90'
$ws.Range("E19").Value = '  +0.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '227.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.26%  '

$ws.Range("E22").Value = '  +0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.378'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1457'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.766'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.503'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.263'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.145'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.207'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05134'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7684'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.834'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.131'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.676'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.68%  '

$ws.Range("D38").Value = '1.302.76'
$ws.Range("E38").Value = '  +1.41%  '

$ws.Range("E39").Value = '  -1.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.698'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9346'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.60%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.815'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.18%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.16%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.613'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.67%  '

$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.983.02'
$ws.Range("E47").Value = '  +0.43%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5198'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.69%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000123'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.767'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07477'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +21.99%  '
